# Update countries & provincias Spain
# - Refresh the "last updated" timestamp in A1 (13:28 -> 14:45)
# - Refresh covid-19 counters for several countries (B:E, G:H columns)
# - Croacia overtakes Tayikistan (both now tied at 8311 total cases, so
#   Croacia - whose count rose - sorts above the unchanged Tayikistan row)
# - Timor Oriental / Santa Lucia swap order (tied counts, alphabetical-ish
#   re-sort) with no numeric change in either row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: refresh timestamp
$ws.Range("A1").Value = "Datos actualizados a 24 de Agosto de 2020 a las 14:45"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5875254
$ws.Range("C4").Value = 1108
$ws.Range("D4").Value = 3167232
$ws.Range("E4").Value = 2527408
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 180614

# Row 17: Arabia Saudita
$ws.Range("B17").Value = 308654
$ws.Range("C17").Value = 1175
$ws.Range("D17").Value = 282888
$ws.Range("E17").Value = 22075
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = 3691

# Row 46: Paises Bajos
$ws.Range("B46").Value = 67128
$ws.Range("C46").Value = 574
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = 6202

# Row 50: Portugal
$ws.Range("B50").Value = 55720
$ws.Range("C50").Value = 123
$ws.Range("D50").Value = 40880
$ws.Range("E50").Value = 13039
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 1801

# Row 55: Ghana
$ws.Range("B55").Value = 43622
$ws.Range("C55").Value = 117
$ws.Range("D55").Value = 41695
$ws.Range("E55").Value = 1664
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 263

# Row 77: Bosnia y Herzegovina
$ws.Range("B77").Value = 18029
$ws.Range("C77").Value = 314
$ws.Range("D77").Value = 11861
$ws.Range("E77").Value = 5621
$ws.Range("G77").Value = 15
$ws.Range("H77").Value = 547

# Row 80: Dinamarca
$ws.Range("B80").Value = 16397
$ws.Range("C80").Value = 80
$ws.Range("D80").Value = 14310
$ws.Range("E80").Value = 1464
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 623

# Row 83: Republica de Macedonia
$ws.Range("B83").Value = 13673
$ws.Range("C83").Value = 78
$ws.Range("D83").Value = 10150
$ws.Range("E83").Value = 2955
$ws.Range("G83").Value = 4
$ws.Range("H83").Value = 568

# Row 98: was Tayikistan -> now Croacia (updated counters)
$ws.Range("A98").Value = "Croacia"
$ws.Range("C98").Value = 136
$ws.Range("D98").Value = 5926
$ws.Range("E98").Value = 2212
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 173

# Row 99: was Croacia -> now Tayikistan (keeps the old, unchanged counters)
$ws.Range("A99").Value = "Tayikistan"
$ws.Range("B99").Value = 8311
$ws.Range("D99").Value = 7108
$ws.Range("E99").Value = 1137
$ws.Range("H99").Value = 66

# Row 101: Finlandia
$ws.Range("E101").Value = 503
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 335

# Row 139: Islandia
$ws.Range("B139").Value = 2073
$ws.Range("C139").Value = 9
$ws.Range("D139").Value = 1946
$ws.Range("E139").Value = 117

# Row 152: Burkina Faso
$ws.Range("B152").Value = 1328
$ws.Range("C152").Value = 8
$ws.Range("D152").Value = 1050
$ws.Range("E152").Value = 223

# Row 202: was Santa Lucia -> now Timor Oriental (counters unchanged)
$ws.Range("A202").Value = "Timor Oriental"

# Row 203: was Timor Oriental -> now Santa Lucia (counters unchanged)
$ws.Range("A203").Value = "Santa Lucia"

# Row 208: Dominica
$ws.Range("B208").Value = 20
$ws.Range("C208").Value = 1
$ws.Range("E208").Value = 2
